# Update the StructureDefinition workbook for the "2025 august" output refresh.
#
# 1) Metadata sheet: the canonical URL moved from the old GitHub-based shorthand
#    repo URL to the new 2rdoc.pt IG URL, and the publication Date was bumped.
# 2) Elements sheet: the column widths were recalculated by the IG Publisher
#    (narrower "best fit" measurements) - apply the new widths.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata sheet (URL + Date) ---------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/mindfulness-notification-enabled"
$meta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# --- 2) Elements sheet (column widths) -------------------------------------
# Excel's ColumnWidth (characters) is quantized to the pixel grid
# (stored_width = round(chars*6)/6 + 5/6), so to land on a target stored
# "width" attribute as closely as possible we invert that formula here
# before assigning ColumnWidth.
$elements = $wb.Worksheets.Item("Elements")

function Set-ColWidthForStoredWidth($ws, [int]$col, [double]$targetWidth) {
    $px = ($targetWidth - (5.0/6.0)) * 6.0
    $pxRound = [Math]::Round($px, 0, [MidpointRounding]::AwayFromZero)
    $chars = $pxRound / 6.0
    $ws.Columns.Item($col).ColumnWidth = $chars
}

Set-ColWidthForStoredWidth $elements 1  16.41796875
Set-ColWidthForStoredWidth $elements 2  16.41796875
Set-ColWidthForStoredWidth $elements 3  9.79296875
Set-ColWidthForStoredWidth $elements 4  7.046875
Set-ColWidthForStoredWidth $elements 5  5.30078125
Set-ColWidthForStoredWidth $elements 6  3.953125
Set-ColWidthForStoredWidth $elements 7  4.265625
Set-ColWidthForStoredWidth $elements 8  12.6875
Set-ColWidthForStoredWidth $elements 9  10.51171875
Set-ColWidthForStoredWidth $elements 11 8.3984375
Set-ColWidthForStoredWidth $elements 15 12.26171875
Set-ColWidthForStoredWidth $elements 20 7.80078125
Set-ColWidthForStoredWidth $elements 21 13.609375
Set-ColWidthForStoredWidth $elements 22 13.91796875
Set-ColWidthForStoredWidth $elements 23 15.01171875
Set-ColWidthForStoredWidth $elements 24 14.62890625
Set-ColWidthForStoredWidth $elements 25 17.08203125
Set-ColWidthForStoredWidth $elements 26 15.18359375
Set-ColWidthForStoredWidth $elements 27 5.07421875
Set-ColWidthForStoredWidth $elements 28 17.98046875
Set-ColWidthForStoredWidth $elements 29 34.578125
Set-ColWidthForStoredWidth $elements 30 13.54296875
Set-ColWidthForStoredWidth $elements 31 11.3203125
Set-ColWidthForStoredWidth $elements 32 15.046875
Set-ColWidthForStoredWidth $elements 33 8.22265625
Set-ColWidthForStoredWidth $elements 34 8.53125
Set-ColWidthForStoredWidth $elements 37 19.5625
